# Apply the edits described by the commit:
#  1. Three tables (on slides 14, 15, 16) switch from the deck's custom
#     "Table_0" style to PowerPoint's built-in default table style.
#  2. The Slide Master theme's colour scheme is swapped from the
#     "Red Violet" / Integral palette to the standard Office palette
#     (i.e. theme1.xml and theme2.xml trade their colour schemes).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{D3DC8939-28C3-4AF2-B50A-4603A5425BE4}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}

# --- 2. Theme colour scheme ------------------------------------------
# Slide master (theme1.xml) goes from the "Red Violet" scheme to the
# standard "Office" scheme. The twelve theme colour slots are reachable
# through the legacy 8-slot ColorScheme collection, which this host
# extends to all twelve theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) via RGB (stored as VBA-style BGR-packed integers).

$cs = $p.SlideMaster.ColorScheme
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
